# Insert a new weekly price record as row 52 ("Fruta / hortaliza, semanal"),
# pushing the existing rows 52-71 down to 53-72.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(52).Insert()

$ws.Cells.Item(52, 1).Value = 3
$ws.Cells.Item(52, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(52, 3).Value = "Coquimbo"
$ws.Cells.Item(52, 4).Value = 45202
$ws.Cells.Item(52, 5).Value = 5
$ws.Cells.Item(52, 6).Value = 300000000
$ws.Cells.Item(52, 7).Value = "Espárragos"
$ws.Cells.Item(52, 8).Value = "Verde"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 900
$ws.Cells.Item(52, 11).Value = 1900
$ws.Cells.Item(52, 12).Value = 1900
$ws.Cells.Item(52, 13).Value = 1900
$ws.Cells.Item(52, 14).Value = "$/kilo"
$ws.Cells.Item(52, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(52, 16).Value = 1900
$ws.Cells.Item(52, 17).Value = 1
$ws.Cells.Item(52, 18).Value = "Hortaliza"
